# Finish edit photon. Update quotes
#
# Adds the "perfect photon from steve" (L/M/N/O, rows 14-17) data block and
# its companion JER-style SQRT formulas (Q/R:T, rows 14-16), mirroring the
# existing A/B:D + G/H:J block. Also updates the sheet view (scroll
# position / active cell) to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New data block: columns L,M,N,O, rows 14-17 ---------------------------
# Row 14 ("y2" perfect-photon-from-steve values)
$ws.Range("L14").Value = 3.49641
$ws.Range("M14").Value = 2.72426
$ws.Range("N14").Value = 2.61667
$ws.Range("O14").Value = 2.7686

# Row 15 (first MUON y2 values)
$ws.Range("L15").Value = 3.56892
$ws.Range("M15").Value = 2.85493
$ws.Range("N15").Value = 2.90771
$ws.Range("O15").Value = 3.08924

# Row 16 (new-photon-with-merging y1 values)
$ws.Range("L16").Value = 3.76354
$ws.Range("M16").Value = 2.8844
$ws.Range("N16").Value = 2.77463
$ws.Range("O16").Value = 2.89704

# Row 17 (second MUON y2 values)
$ws.Range("L17").Value = 3.76815
$ws.Range("M17").Value = 2.90265
$ws.Range("N17").Value = 2.86575
$ws.Range("O17").Value = 2.92195

# --- New formulas: columns Q (standalone) and R:T (shared), rows 14-16 -----
$ws.Range("Q14").Formula = "=SQRT(L15^2-L14^2)"
$ws.Range("R14:T14").Formula = "=SQRT(M15^2-M14^2)"

$ws.Range("Q15").Formula = "=SQRT(L16^2-L14^2)"
$ws.Range("R15:T15").Formula = "=SQRT(M16^2-M14^2)"

$ws.Range("Q16").Formula = "=SQRT(L17^2-L14^2)"
$ws.Range("R16:T16").Formula = "=SQRT(M17^2-M14^2)"

# --- Sheet view: move the selection (also clears the old topLeftCell) -----
$ws.Range("I23").Select()
